# Update "想去人数" (interested-count) figures that changed between scrapes.
$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (Exhibitions)
$wsExpo = $wb.Worksheets.Item(1)
$wsExpo.Range("F2").Value = 82
$wsExpo.Range("F13").Value = 1477
$wsExpo.Range("F15").Value = 2735

# Sheet 4: 全部类型 (All types) mirrors the same events
$wsAll = $wb.Worksheets.Item(4)
$wsAll.Range("F2").Value = 82
$wsAll.Range("F16").Value = 1477
$wsAll.Range("F18").Value = 2735
